$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 19-21: the IndianWings (Bangalore) row moves from row 19 down
# to row 21, while the two EuropeanWings rows (previously rows 20 and 21)
# shift up to rows 19 and 20.
$row19 = @("EuropeanWings", "Paris-Beauvais-Tille", "LFOB", "Hungary-Budapest-Listz", "LHBP")
$row20 = @("EuropeanWings", "Hungary-Budapest-Listz", "LHBP", "Paris-Beauvais-Tille", "LFOB")
$row21 = @("IndianWings", "Bangalore-India", "VOBL", "Indira Gandhi Intl New-Delhi-India", "VIDP")
$row22 = @("IndianWings", "Mumbai-India", "VABB", "Calcutta-India", "VECC")
$row23 = @("IndianWings", "Chennai-India", "VOMM", "Jaipur-India", "VIJP")

# Clear any pre-existing formatting on rows 19-23 so old style placement
# doesn't bleed through once the rows are rewritten with new content.
$ws.Range("A19:E23").ClearFormats()

for ($c = 1; $c -le 5; $c++) {
    $ws.Cells.Item(19, $c).Value = $row19[$c - 1]
    $ws.Cells.Item(20, $c).Value = $row20[$c - 1]
    $ws.Cells.Item(21, $c).Value = $row21[$c - 1]
    $ws.Cells.Item(22, $c).Value = $row22[$c - 1]
    $ws.Cells.Item(23, $c).Value = $row23[$c - 1]
}

# Match vertical-center alignment style exactly as it lands per-row (the
# style travelled with the row's original cell position when the rows were
# reordered, so row 20 carries it on column B instead of D).
$ws.Range("D19").VerticalAlignment = -4108
$ws.Range("B20").VerticalAlignment = -4108
$ws.Range("D21").VerticalAlignment = -4108
$ws.Range("D22").VerticalAlignment = -4108
$ws.Range("D23").VerticalAlignment = -4108

$ws.Range("E23").Select()
